$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rowsForView")

# Update letter title text in D25 / D26 (shared string)
$ws.Range("D25").Value = "Lettre de relance - dernier rappel"
$ws.Range("D26").Value = "Lettre de relance - dernier rappel"

# Update last-updated date in F25 / F26 from 2021-04-21 to 2021-04-28
$ws.Range("F25").Value = [DateTime]"2021-04-28"
$ws.Range("F26").Value = [DateTime]"2021-04-28"

# Update window size recorded in the workbook view
$excel.ActiveWindow.Width = 20415
$excel.ActiveWindow.Height = 6990
